$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - time ranges
$ws.Range("E6").Value = "9->12"
$ws.Range("F6").Value = "9->11"
$ws.Range("G6").Value = "9->12"

# Row 7 - time ranges
$ws.Range("D7").Value = "6->7"
$ws.Range("E7").Value = "2->5"
$ws.Range("F7").Value = "11->5"

# Row 8 - new meeting label
$ws.Range("F8").Value = "Meeting 10"

# Row 9 - updated note text + new meeting note
$ws.Range("C9").Value = "Meet with TAS about whats going on for the week"
$ws.Range("F9").Value = "Met with John, Lorna, Craig about direction of the project"

# Row 11 - meeting notes text (now holds the "Talked with john..." note)
$ws.Range("C11").Value = "Talked with john about options going forward. `n- Search capability`n- How many australian calls`n- Time trials`n- $$ Cost to map goata`n- What do they need in order to map for us"

# Selection moved to F9
$ws.Range("F9").Select()
